# "Update Data Sources from LFX" automated refresh
# --------------------------------------------------
# 1) Six data-source tables get re-pointed from the old "Table_0" table
#    style to a new built-in table style.
# 2) The deck's theme-color palette is refreshed (the LFX bot swaps the
#    "LF Energy Theme 2023"/Geometric palette and the "Simple Light"
#    palette between the two themes shipped in the template).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Re-style every table in the deck with the new table-style GUID.
# ---------------------------------------------------------------------
$newTableStyleId = "{8FA931E9-57A9-4F2B-8CB2-13A407556EF8}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2. Refresh the theme colour scheme to the "LF Energy Theme 2023"
#    (Geometric) palette.
# ---------------------------------------------------------------------
function Convert-HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme's 1-based indices:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$geometricPalette = @(
    "222222",
    "FFFFFF",
    "434343",
    "999999",
    "003778",
    "0094FF",
    "5B1DE7",
    "12E2E2",
    "FF00AA",
    "ACDE1F",
    "0077CC",
    "F06292"
)

$design = $p.Designs.Item($p.Designs.Count)
$colorScheme = $design.SlideMaster.Theme.ThemeColorScheme

for ($k = 1; $k -le $colorScheme.Count; $k++) {
    $colorScheme.Item($k).RGB = Convert-HexToRgb $geometricPalette[$k - 1]
}
